# Regenerate s_vals data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP) and G (sum) for rows 2-16.
# Column A (date) and column F (Win) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    # row, B(TB), C(d2S), D(K), E(IP), G(sum)
    @(2,  0.7287194209349384, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 3.594575437922795),
    @(3,  1.505614041169197,  1.65323645889881,  0.1529057820181812, 0.4998867070740569, 3.811642989160245),
    @(4,  0.7287194209349384, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 3.594575437922795),
    @(5,  1.505614041169197,  1.65323645889881,  0.1529057820181812, 6.48142807727062,   9.793184359356808),
    @(6,  3.182878228561681,  1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538),
    @(7,  3.182878228561681,  1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538),
    @(8,  0.7287194209349384, 0.3375848360084654,0.1529057820181812, 0.4998867070740569, 1.719096746035642),
    @(9,  0.7287194209349384, 0.3375848360084654,0.1529057820181812, 0.4998867070740569, 1.719096746035642),
    @(10, 1.505614041169197,  86.29678392075563, 10137753.70137369,  6.48142807727062,   10137847.98519973),
    @(11, 3.182878228561681,  1.65323645889881,  0.1529057820181812, 0.4998867070740569, 5.488907176552729),
    @(12, 1.505614041169197,  1.65323645889881,  0.7127328510149897, 6.48142807727062,   10.35301142835362),
    @(13, 0.006876353814593728,0.05231270169004087,0.1529057820181812,246.9852506941017, 247.1973455316245),
    @(14, 0.7287194209349384, 1.65323645889881,  0.7127328510149897, 6.48142807727062,   9.576116808119359),
    @(15, 3.182878228561681,  1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538),
    @(16, 1.505614041169197,  1.65323645889881,  0.1529057820181812, 0.4998867070740569, 3.811642989160245)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]   # B: TB
    $ws.Cells.Item($row, 3).Value = $entry[2]   # C: d2S
    $ws.Cells.Item($row, 4).Value = $entry[3]   # D: K
    $ws.Cells.Item($row, 5).Value = $entry[4]   # E: IP
    $ws.Cells.Item($row, 7).Value = $entry[5]   # G: sum
}
